$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.997.50'
$ws.Cells.Item(2, 5).Value = '  -0.23%  '
$ws.Cells.Item(3, 4).Value = '2.526.42'
$ws.Cells.Item(3, 5).Value = '  +0.31%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '536.11'
$ws.Cells.Item(5, 5).Value = '  -0.07%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '137.74'
$ws.Cells.Item(6, 5).Value = '  -1.48%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.998'
$ws.Cells.Item(7, 5).Value = '  -0.16%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.567'
$ws.Cells.Item(8, 5).Value = '  +1.09%  '
$ws.Cells.Item(9, 4).Value = '2.527.39'
$ws.Cells.Item(9, 5).Value = '  +0.29%  '
$ws.Cells.Item(10, 5).Value = '  +1.92%  '
$ws.Cells.Item(11, 5).Value = '  -0.22%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '5.39'
$ws.Cells.Item(12, 5).Value = '  -0.64%  '
$ws.Cells.Item(13, 5).Value = '  -2.01%  '
$ws.Cells.Item(14, 4).Value = '2.948.35'
$ws.Cells.Item(14, 5).Value = '  -0.56%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '23.27'
$ws.Cells.Item(15, 5).Value = '  +1.58%  '
$ws.Cells.Item(16, 4).Value = '59.147.44'
$ws.Cells.Item(16, 5).Value = '  +0.08%  '
$ws.Cells.Item(17, 5).Value = '  -0.06%  '
$ws.Cells.Item(18, 4).Value = '2.513.11'
$ws.Cells.Item(18, 5).Value = '  -0.39%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.09'
$ws.Cells.Item(19, 5).Value = '  +1.81%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '4.27'
$ws.Cells.Item(20, 5).Value = '  +1.21%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '325.74'
$ws.Cells.Item(21, 5).Value = '  +1.34%  '
$ws.Cells.Item(22, 5).Value = '  +0.46%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.88'
$ws.Cells.Item(23, 5).Value = '  +1.02%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '64.93'
$ws.Cells.Item(24, 5).Value = '  +4.02%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.421'
$ws.Cells.Item(25, 5).Value = '  -0.23%  '
$ws.Cells.Item(26, 5).Value = '  +0.49%  '
$ws.Cells.Item(27, 5).Value = '  +0.85%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.63'
$ws.Cells.Item(28, 5).Value = '  -1.52%  '
$ws.Cells.Item(29, 2).Value = 'PEPE'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(29, 4).Value = '0.0₃0779'
$ws.Cells.Item(29, 5).Value = '  +1.82%  '
$ws.Cells.Item(30, 2).Value = 'Aptos'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.74'
$ws.Cells.Item(30, 5).Value = '  +0.36%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.78'
$ws.Cells.Item(31, 5).Value = '  -1.11%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '167.35'
$ws.Cells.Item(32, 5).Value = '  +4.11%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.18'
$ws.Cells.Item(33, 5).Value = '  +4.63%  '
$ws.Cells.Item(34, 5).Value = '  -0.12%  '
$ws.Cells.Item(35, 5).Value = '  -3.34%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '18.60'
$ws.Cells.Item(36, 5).Value = '  +0.55%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.14'
$ws.Cells.Item(37, 5).Value = '  -1.41%  '
$ws.Cells.Item(38, 5).Value = '  -0.32%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '36.73'
$ws.Cells.Item(39, 5).Value = '  -0.57%  '
$ws.Cells.Item(40, 5).Value = '  +3.99%  '
$ws.Cells.Item(41, 5).Value = '  +0.31%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '5.24'
$ws.Cells.Item(42, 5).Value = '  -0.22%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '282.15'
$ws.Cells.Item(43, 5).Value = '  -0.15%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.997'
$ws.Cells.Item(44, 5).Value = '  -0.06%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.606'
$ws.Cells.Item(45, 5).Value = '  +1.93%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '10.87'
$ws.Cells.Item(46, 5).Value = '  -0.07%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '128.28'
$ws.Cells.Item(47, 5).Value = '  +4.91%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0513'
$ws.Cells.Item(49, 5).Value = '  +0.66%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0224'
$ws.Cells.Item(50, 5).Value = '  +0.59%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '17.40'
$ws.Cells.Item(51, 5).Value = '  +0.13%  '
